{"js": "// Merge the split \"<id>...</id>\" runs back into a single run for the\n// p033r_1 and p033r_2 paragraphs (newly downloaded tc/tcn/tl content was\n// re-joined into one text run, keeping the first run's formatting).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The paragraphs that must be collapsed to a single run are exactly the\n// ones whose *entire* text is \"<id>p033r_1</id>\" or \"<id>p033r_2</id>\"\n// (the sibling \"<id>fig_p033r_*</id>\" paragraphs are left untouched).\nconst targetTexts = new Set([\"<id>p033r_1</id>\", \"<id>p033r_2</id>\"]);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (targetTexts.has(paragraph.text)) {\n    // Re-inserting the same text over the whole paragraph range collapses\n    // the existing multi-run split into a single run that carries the\n    // formatting of the range's first run (the Courier-New/brown id tag\n    // styling), exactly matching the merged run in the target XML.\n    const range = paragraph.getRange();\n    range.insertText(paragraph.text, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Merge the split \"<id>...</id>\" runs back into a single run for the\n# p033r_1 and p033r_2 paragraphs (newly downloaded tc/tcn/tl content was\n# re-joined into one text run, keeping the first run's formatting).\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$d = $word.ActiveDocument\n\n$targets = @(\"<id>p033r_1</id>\", \"<id>p033r_2</id>\")\n\nforeach ($targetText in $targets) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # Re-\"replacing\" the found text with itself collapses the existing\n    # multi-run split (separate \"<id>\", \"p033rN_\", \"</id>\" runs) into a\n    # single run that carries the formatting of the match's first run\n    # (the Courier-New / brown id-tag styling), exactly matching the\n    # merged run in the target document.\n    $d.Content.Find.Execute($targetText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $targetText, $wdReplaceOne) | Out-Null\n}\n"}
